$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(3, 1).Value = 99.69248304756557
$ws.Cells.Item(4, 1).Value = 99.65015726742229
$ws.Cells.Item(5, 1).Value = 99.94069548291762
$ws.Cells.Item(6, 1).Value = 100.5054026859905
$ws.Cells.Item(7, 1).Value = 100.4480682106593
$ws.Cells.Item(8, 1).Value = 100.3492894208657
$ws.Cells.Item(9, 1).Value = 100.407706825512
$ws.Cells.Item(10, 1).Value = 100.2094256881611
$ws.Cells.Item(11, 1).Value = 100.7184991547432
$ws.Cells.Item(12, 1).Value = 100.6269414885737
$ws.Cells.Item(13, 1).Value = 100.647015540968
$ws.Cells.Item(14, 1).Value = 100.6792753401545
$ws.Cells.Item(15, 1).Value = 100.5069399081758
$ws.Cells.Item(16, 1).Value = 100.7657247984114
$ws.Cells.Item(17, 1).Value = 101.0086459722758
$ws.Cells.Item(18, 1).Value = 101.4555688854043
$ws.Cells.Item(19, 1).Value = 101.5108438758932
$ws.Cells.Item(20, 1).Value = 101.6599133315943
$ws.Cells.Item(21, 1).Value = 101.5957530987969
$ws.Cells.Item(22, 1).Value = 101.3308133892882
$ws.Cells.Item(23, 1).Value = 101.3426609138731
$ws.Cells.Item(24, 1).Value = 101.7229995607766
$ws.Cells.Item(25, 1).Value = 101.9198894098388
$ws.Cells.Item(26, 1).Value = 102.0251155140434
$ws.Cells.Item(27, 1).Value = 102.1294649626944
$ws.Cells.Item(28, 1).Value = 102.5761973162791
$ws.Cells.Item(29, 1).Value = 102.2536748638381
$ws.Cells.Item(30, 1).Value = 102.1451134991974
$ws.Cells.Item(31, 1).Value = 101.6282350723164
$ws.Cells.Item(32, 1).Value = 102.0376878779384
$ws.Cells.Item(33, 1).Value = 102.3233665619997
$ws.Cells.Item(34, 1).Value = 102.3824525958246
$ws.Cells.Item(35, 1).Value = 102.5578007949158
$ws.Cells.Item(36, 1).Value = 102.8306754933793
$ws.Cells.Item(37, 1).Value = 102.9043236480488
$ws.Cells.Item(38, 1).Value = 102.6928712953071
$ws.Cells.Item(39, 1).Value = 103.2039150963427
$ws.Cells.Item(40, 1).Value = 103.0254962151673
$ws.Cells.Item(41, 1).Value = 103.3323560149197
$ws.Cells.Item(42, 1).Value = 103.5273386636665
$ws.Cells.Item(43, 1).Value = 103.4921387825334
$ws.Cells.Item(44, 1).Value = 103.0819582219913
$ws.Cells.Item(45, 1).Value = 103.1013248976
$ws.Cells.Item(46, 1).Value = 103.206055602845
$ws.Cells.Item(47, 1).Value = 103.2801165510631
$ws.Cells.Item(48, 1).Value = 102.9489430993215
$ws.Cells.Item(49, 1).Value = 102.785469726719
$ws.Cells.Item(50, 1).Value = 102.5490276553974
$ws.Cells.Item(51, 1).Value = 102.615654804477
$ws.Cells.Item(52, 1).Value = 102.7663487994706
$ws.Cells.Item(53, 1).Value = 102.6837178077704
$ws.Cells.Item(54, 1).Value = 102.2880694682678
$ws.Cells.Item(55, 1).Value = 102.8259125989839
$ws.Cells.Item(56, 1).Value = 102.9741234359998
$ws.Cells.Item(57, 1).Value = 102.9407992119655
$ws.Cells.Item(58, 1).Value = 103.0600738975572
$ws.Cells.Item(59, 1).Value = 103.1361116848377
$ws.Cells.Item(60, 1).Value = 102.9341927618848
$ws.Cells.Item(61, 1).Value = 103.1878374975335
$ws.Cells.Item(62, 1).Value = 102.9326865863377
$ws.Cells.Item(63, 1).Value = 103.2512559347332
$ws.Cells.Item(64, 1).Value = 103.2232025962558
$ws.Cells.Item(65, 1).Value = 103.0654852714342
$ws.Cells.Item(66, 1).Value = 103.1788689021445
$ws.Cells.Item(67, 1).Value = 103.4282119748315
$ws.Cells.Item(68, 1).Value = 103.2571071025947
$ws.Cells.Item(69, 1).Value = 103.4283884469656
$ws.Cells.Item(70, 1).Value = 103.7308324570732
$ws.Cells.Item(71, 1).Value = 103.1001321066149
$ws.Cells.Item(72, 1).Value = 103.2583639536984
$ws.Cells.Item(73, 1).Value = 103.4365467971583
$ws.Cells.Item(74, 1).Value = 103.9778504360769
$ws.Cells.Item(75, 1).Value = 103.6040955408012
$ws.Cells.Item(76, 1).Value = 103.3526811524504
$ws.Cells.Item(77, 1).Value = 103.3523004893574
$ws.Cells.Item(78, 1).Value = 103.2131864244591
$ws.Cells.Item(79, 1).Value = 103.4274770734536
$ws.Cells.Item(80, 1).Value = 103.9991400559845
$ws.Cells.Item(81, 1).Value = 104.024780250218
$ws.Cells.Item(82, 1).Value = 103.7886003362273
$ws.Cells.Item(83, 1).Value = 104.3559236010447
$ws.Cells.Item(84, 1).Value = 104.4866599665155
$ws.Cells.Item(85, 1).Value = 103.7821582176173
$ws.Cells.Item(86, 1).Value = 103.9354844224032
$ws.Cells.Item(87, 1).Value = 103.9425652649151
$ws.Cells.Item(88, 1).Value = 104.0945858642456
$ws.Cells.Item(89, 1).Value = 103.8245338066157
$ws.Cells.Item(90, 1).Value = 103.1854564187146
$ws.Cells.Item(91, 1).Value = 103.101327738135
$ws.Cells.Item(92, 1).Value = 103.1509158234016
$ws.Cells.Item(93, 1).Value = 102.9140793576466
$ws.Cells.Item(94, 1).Value = 103.5487486874957
$ws.Cells.Item(95, 1).Value = 103.7170642390891
$ws.Cells.Item(96, 1).Value = 103.5374317463395
$ws.Cells.Item(97, 1).Value = 104.0858795098562
$ws.Cells.Item(98, 1).Value = 104.2198742582753
$ws.Cells.Item(99, 1).Value = 104.2729838348916
$ws.Cells.Item(100, 1).Value = 104.6002058209827
$ws.Cells.Item(101, 1).Value = 104.7843860092436
$ws.Cells.Item(102, 1).Value = 105.1941207871387
$ws.Cells.Item(103, 1).Value = 104.4839750267953
$ws.Cells.Item(104, 1).Value = 104.2927045241435
$ws.Cells.Item(105, 1).Value = 104.3196418545777
$ws.Cells.Item(106, 1).Value = 104.6147733552671
$ws.Cells.Item(107, 1).Value = 104.7001317890918
$ws.Cells.Item(108, 1).Value = 104.7325956486104
$ws.Cells.Item(109, 1).Value = 104.2410116914176
$ws.Cells.Item(110, 1).Value = 104.4546603389535
$ws.Cells.Item(111, 1).Value = 104.3114337980369
$ws.Cells.Item(112, 1).Value = 104.6087661564664
$ws.Cells.Item(113, 1).Value = 104.3438438072361
$ws.Cells.Item(114, 1).Value = 104.4862144872693
$ws.Cells.Item(115, 1).Value = 103.7245249305959
$ws.Cells.Item(116, 1).Value = 103.734254818643
$ws.Cells.Item(117, 1).Value = 103.8857156634132
$ws.Cells.Item(118, 1).Value = 103.8001096003688
$ws.Cells.Item(119, 1).Value = 103.8102678934546
$ws.Cells.Item(120, 1).Value = 103.8351534659642
$ws.Cells.Item(121, 1).Value = 104.0323829099598
$ws.Cells.Item(122, 1).Value = 104.0888057584971
$ws.Cells.Item(123, 1).Value = 104.2588571961974
$ws.Cells.Item(124, 1).Value = 104.5029658584868
$ws.Cells.Item(125, 1).Value = 104.3551284039529
$ws.Cells.Item(126, 1).Value = 104.3551284039529
$ws.Cells.Item(127, 1).Value = 104.2723977789433
$ws.Cells.Item(128, 1).Value = 104.2186993565057
$ws.Cells.Item(129, 1).Value = 104.4392333527118
$ws.Cells.Item(130, 1).Value = 104.3376542589845
$ws.Cells.Item(131, 1).Value = 104.571978832556
$ws.Cells.Item(132, 1).Value = 104.9334909407579
$ws.Cells.Item(133, 1).Value = 104.9768323857903
$ws.Cells.Item(134, 1).Value = 105.0790535009185
$ws.Cells.Item(135, 1).Value = 104.9692558754699
$ws.Cells.Item(136, 1).Value = 105.6825546483224
$ws.Cells.Item(137, 1).Value = 105.7004043614803
$ws.Cells.Item(138, 1).Value = 105.8830636817305
$ws.Cells.Item(139, 1).Value = 105.7800252937052
$ws.Cells.Item(140, 1).Value = 105.6647999206322
$ws.Cells.Item(141, 1).Value = 105.8129261273272
$ws.Cells.Item(142, 1).Value = 105.9391564610592
$ws.Cells.Item(143, 1).Value = 105.9825971896758
$ws.Cells.Item(144, 1).Value = 105.7154320831238
$ws.Cells.Item(145, 1).Value = 105.2401979724218
$ws.Cells.Item(146, 1).Value = 105.0779647045263
$ws.Cells.Item(147, 1).Value = 104.4392997294353
$ws.Cells.Item(148, 1).Value = 104.7474075617528
$ws.Cells.Item(149, 1).Value = 104.9613883852087
$ws.Cells.Item(150, 1).Value = 105.4543048941611
$ws.Cells.Item(151, 1).Value = 104.8066183477122
$ws.Cells.Item(152, 1).Value = 104.8215956892161
$ws.Cells.Item(153, 1).Value = 104.820866706809
$ws.Cells.Item(154, 1).Value = 104.8159920396024
$ws.Cells.Item(155, 1).Value = 105.454616140817
$ws.Cells.Item(156, 1).Value = 105.679867964786
$ws.Cells.Item(157, 1).Value = 105.1867334761358
$ws.Cells.Item(158, 1).Value = 104.205225915887
$ws.Cells.Item(159, 1).Value = 104.3785319990359
$ws.Cells.Item(160, 1).Value = 103.8931790504332
$ws.Cells.Item(161, 1).Value = 104.1825987265084
$ws.Cells.Item(162, 1).Value = 103.6682509611537
$ws.Cells.Item(163, 1).Value = 103.5207639705611
$ws.Cells.Item(164, 1).Value = 103.6052181813543
$ws.Cells.Item(165, 1).Value = 103.4109581621887
$ws.Cells.Item(166, 1).Value = 103.1334554917413
$ws.Cells.Item(167, 1).Value = 102.212958325392
$ws.Cells.Item(168, 1).Value = 102.4067515533737
$ws.Cells.Item(169, 1).Value = 102.467117583097
$ws.Cells.Item(170, 1).Value = 102.4326788778613
$ws.Cells.Item(171, 1).Value = 102.2878771893793
$ws.Cells.Item(172, 1).Value = 103.0005874123444
$ws.Cells.Item(173, 1).Value = 103.2229958339958
$ws.Cells.Item(174, 1).Value = 102.8074024879286
$ws.Cells.Item(175, 1).Value = 102.6796004911194
$ws.Cells.Item(176, 1).Value = 102.978767806976
$ws.Cells.Item(177, 1).Value = 103.1251699630746
$ws.Cells.Item(178, 1).Value = 103.0880983555124
$ws.Cells.Item(179, 1).Value = 102.8921019982229
$ws.Cells.Item(180, 1).Value = 103.0688065472065
$ws.Cells.Item(181, 1).Value = 103.5766399556518
$ws.Cells.Item(182, 1).Value = 103.1916567840573
$ws.Cells.Item(183, 1).Value = 102.7085691477715
$ws.Cells.Item(184, 1).Value = 102.3537231148896
$ws.Cells.Item(185, 1).Value = 102.4286218730058
$ws.Cells.Item(186, 1).Value = 102.9676204729249
$ws.Cells.Item(187, 1).Value = 102.4292647236852
$ws.Cells.Item(188, 1).Value = 102.8653272247564
$ws.Cells.Item(189, 1).Value = 102.8693159072993
$ws.Cells.Item(190, 1).Value = 103.3156614670835
$ws.Cells.Item(191, 1).Value = 102.8791668812905
$ws.Cells.Item(192, 1).Value = 103.2984944434431
$ws.Cells.Item(193, 1).Value = 103.1354293926127
$ws.Cells.Item(194, 1).Value = 103.4025844294447
$ws.Cells.Item(195, 1).Value = 103.4083265795245
$ws.Cells.Item(196, 1).Value = 103.4254027572152
$ws.Cells.Item(197, 1).Value = 103.8493440962932
$ws.Cells.Item(198, 1).Value = 103.8759856298659
$ws.Cells.Item(199, 1).Value = 103.4499157061046
$ws.Cells.Item(200, 1).Value = 103.5814791225
$ws.Cells.Item(201, 1).Value = 103.2199644614203
$ws.Cells.Item(202, 1).Value = 103.4584357914995
$ws.Cells.Item(203, 1).Value = 102.9483419657644
$ws.Cells.Item(204, 1).Value = 103.0314494789732
$ws.Cells.Item(205, 1).Value = 103.3910462540462
$ws.Cells.Item(206, 1).Value = 103.7300254717483
$ws.Cells.Item(207, 1).Value = 104.3585419394528
$ws.Cells.Item(208, 1).Value = 104.4221128273882
$ws.Cells.Item(209, 1).Value = 104.3459919996829
$ws.Cells.Item(210, 1).Value = 104.7773801801366
$ws.Cells.Item(211, 1).Value = 104.933884409216
$ws.Cells.Item(212, 1).Value = 104.698206101121
$ws.Cells.Item(213, 1).Value = 105.0036306443514
$ws.Cells.Item(214, 1).Value = 105.0669800395573
$ws.Cells.Item(215, 1).Value = 105.108282957487
$ws.Cells.Item(216, 1).Value = 105.6763410711545
$ws.Cells.Item(217, 1).Value = 105.887445342585
$ws.Cells.Item(218, 1).Value = 105.2502297324288
$ws.Cells.Item(219, 1).Value = 105.4753614647813
$ws.Cells.Item(220, 1).Value = 105.6178390604507
$ws.Cells.Item(221, 1).Value = 105.4017999301665
$ws.Cells.Item(222, 1).Value = 105.5358440303272
$ws.Cells.Item(223, 1).Value = 105.1058743573775
$ws.Cells.Item(224, 1).Value = 105.0445702613885
$ws.Cells.Item(225, 1).Value = 105.0819645137914
$ws.Cells.Item(226, 1).Value = 105.0939747350474
$ws.Cells.Item(227, 1).Value = 105.3630102317705
$ws.Cells.Item(228, 1).Value = 105.5349943684972
$ws.Cells.Item(229, 1).Value = 105.4222131141887
$ws.Cells.Item(230, 1).Value = 105.5432104252085
$ws.Cells.Item(231, 1).Value = 106.0152797039407
$ws.Cells.Item(232, 1).Value = 105.5551319471562
$ws.Cells.Item(233, 1).Value = 105.3909685399797
$ws.Cells.Item(234, 1).Value = 105.4347295032984
$ws.Cells.Item(235, 1).Value = 105.5925040692154
$ws.Cells.Item(236, 1).Value = 105.5439996209943
$ws.Cells.Item(237, 1).Value = 105.4692385714453
$ws.Cells.Item(238, 1).Value = 104.9661488442551
$ws.Cells.Item(239, 1).Value = 105.1106289868407
$ws.Cells.Item(240, 1).Value = 104.8962854860906
$ws.Cells.Item(241, 1).Value = 105.1563947301561
$ws.Cells.Item(242, 1).Value = 105.4802018137399
$ws.Cells.Item(243, 1).Value = 105.5221863448165
$ws.Cells.Item(244, 1).Value = 105.4665796674621
$ws.Cells.Item(245, 1).Value = 105.5665191924277
$ws.Cells.Item(246, 1).Value = 105.6789134973537
$ws.Cells.Item(247, 1).Value = 105.7509011899669
$ws.Cells.Item(248, 1).Value = 105.4450785792067
$ws.Cells.Item(249, 1).Value = 105.7947025665379
$ws.Cells.Item(250, 1).Value = 105.7947025665379
$ws.Cells.Item(251, 1).Value = 106.3215956736972
$ws.Cells.Item(252, 1).Value = 106.1184309843087
$ws.Cells.Item(253, 1).Value = 106.2149000347801
$ws.Cells.Item(254, 1).Value = 106.3196536207832
$ws.Cells.Item(255, 1).Value = 106.8142506257778
$ws.Cells.Item(256, 1).Value = 106.255660875949
$ws.Cells.Item(257, 1).Value = 105.9069434302788
$ws.Cells.Item(258, 1).Value = 106.1046999554834
$ws.Cells.Item(259, 1).Value = 106.2924012260043
$ws.Cells.Item(260, 1).Value = 105.9678977632355
$ws.Cells.Item(261, 1).Value = 105.4521769212792
$ws.Cells.Item(262, 1).Value = 105.0930165705395
$ws.Cells.Item(263, 1).Value = 105.5059959910945
$ws.Cells.Item(264, 1).Value = 105.4526168340259
$ws.Cells.Item(265, 1).Value = 105.9733802421562
$ws.Cells.Item(266, 1).Value = 106.0182830977236
$ws.Cells.Item(267, 1).Value = 106.045421858034
$ws.Cells.Item(268, 1).Value = 105.7764831010053
$ws.Cells.Item(269, 1).Value = 105.6071096717094
$ws.Cells.Item(270, 1).Value = 105.8071763700755
$ws.Cells.Item(271, 1).Value = 105.7796156179031
$ws.Cells.Item(272, 1).Value = 105.8494450572957
$ws.Cells.Item(273, 1).Value = 105.7869637729499
$ws.Cells.Item(274, 1).Value = 105.51062209336
$ws.Cells.Item(275, 1).Value = 105.6973432008014
$ws.Cells.Item(276, 1).Value = 105.4761306304697
$ws.Cells.Item(277, 1).Value = 105.7431196386288
$ws.Cells.Item(278, 1).Value = 105.7889279165069
$ws.Cells.Item(279, 1).Value = 106.6418311336119
$ws.Cells.Item(280, 1).Value = 106.3478370838757
$ws.Cells.Item(281, 1).Value = 106.2425919359347
$ws.Cells.Item(282, 1).Value = 106.6083249655228
$ws.Cells.Item(283, 1).Value = 106.4817254672631
$ws.Cells.Item(284, 1).Value = 106.3344440409818
$ws.Cells.Item(285, 1).Value = 106.7049782663946
$ws.Cells.Item(286, 1).Value = 106.4934586926806
$ws.Cells.Item(287, 1).Value = 106.3223898785538
$ws.Cells.Item(288, 1).Value = 106.0344919365438
$ws.Cells.Item(289, 1).Value = 105.835764391773
$ws.Cells.Item(290, 1).Value = 105.9756681013528
$ws.Cells.Item(291, 1).Value = 105.9045225014811
$ws.Cells.Item(292, 1).Value = 105.4946454180495
$ws.Cells.Item(293, 1).Value = 105.4088692777462
$ws.Cells.Item(294, 1).Value = 105.5543976128405
$ws.Cells.Item(295, 1).Value = 105.3274328363153
$ws.Cells.Item(296, 1).Value = 105.1474198270827
$ws.Cells.Item(297, 1).Value = 105.1996440114214
$ws.Cells.Item(298, 1).Value = 105.3745366985605
$ws.Cells.Item(299, 1).Value = 105.6730945994259
$ws.Cells.Item(300, 1).Value = 105.5057601726612
$ws.Cells.Item(301, 1).Value = 104.9780219289257
$ws.Cells.Item(302, 1).Value = 104.6310898449762
$ws.Cells.Item(303, 1).Value = 104.132069826613
$ws.Cells.Item(304, 1).Value = 103.8451380794763
$ws.Cells.Item(305, 1).Value = 104.0932344477759
$ws.Cells.Item(306, 1).Value = 103.7802060060735
$ws.Cells.Item(307, 1).Value = 103.5848237364097
$ws.Cells.Item(308, 1).Value = 103.0792404972172
$ws.Cells.Item(309, 1).Value = 102.9478065757367
$ws.Cells.Item(310, 1).Value = 103.2004859678956
$ws.Cells.Item(311, 1).Value = 103.3487527607641
$ws.Cells.Item(312, 1).Value = 103.0267251239777
$ws.Cells.Item(313, 1).Value = 103.6997162108527
$ws.Cells.Item(314, 1).Value = 103.0022901737188
$ws.Cells.Item(315, 1).Value = 103.3772185461961
$ws.Cells.Item(316, 1).Value = 103.8064668541227
$ws.Cells.Item(317, 1).Value = 103.4036922795527
$ws.Cells.Item(318, 1).Value = 103.0199510156762
$ws.Cells.Item(319, 1).Value = 102.6795056241313
$ws.Cells.Item(320, 1).Value = 103.3211687990152
$ws.Cells.Item(321, 1).Value = 103.6176985005455
$ws.Cells.Item(322, 1).Value = 103.2105773596777
$ws.Cells.Item(323, 1).Value = 103.4881943495182
$ws.Cells.Item(324, 1).Value = 103.4402700475564
$ws.Cells.Item(325, 1).Value = 103.4657865147114
$ws.Cells.Item(326, 1).Value = 103.6386133436894
$ws.Cells.Item(327, 1).Value = 103.8559724628568
$ws.Cells.Item(328, 1).Value = 103.5862777105362
$ws.Cells.Item(329, 1).Value = 104.2737875316312
$ws.Cells.Item(330, 1).Value = 104.4610821003216
$ws.Cells.Item(331, 1).Value = 104.1202140394836
$ws.Cells.Item(332, 1).Value = 104.1484517352789
$ws.Cells.Item(333, 1).Value = 104.5104570021654
$ws.Cells.Item(334, 1).Value = 105.0096126541378
$ws.Cells.Item(335, 1).Value = 105.0547575977351
$ws.Cells.Item(336, 1).Value = 105.4503162041613
$ws.Cells.Item(337, 1).Value = 105.8191984814342
$ws.Cells.Item(338, 1).Value = 105.8626894306674
$ws.Cells.Item(339, 1).Value = 106.0740662634622
$ws.Cells.Item(340, 1).Value = 105.888544507271
$ws.Cells.Item(341, 1).Value = 105.5687687967225
$ws.Cells.Item(342, 1).Value = 105.2179119213302
$ws.Cells.Item(343, 1).Value = 105.3548219312665
$ws.Cells.Item(344, 1).Value = 105.3997276713164
$ws.Cells.Item(345, 1).Value = 105.1234910404208
$ws.Cells.Item(346, 1).Value = 105.3718239718268
$ws.Cells.Item(347, 1).Value = 105.1007979767237
$ws.Cells.Item(348, 1).Value = 105.2372117848647
$ws.Cells.Item(349, 1).Value = 105.175644167342
$ws.Cells.Item(350, 1).Value = 105.0683759323845
$ws.Cells.Item(351, 1).Value = 105.388865453885
$ws.Cells.Item(352, 1).Value = 105.9213924093219
$ws.Cells.Item(353, 1).Value = 105.7967384867399
$ws.Cells.Item(354, 1).Value = 105.8935728041645
$ws.Cells.Item(355, 1).Value = 105.5965951048511
$ws.Cells.Item(356, 1).Value = 105.5811909949901
$ws.Cells.Item(357, 1).Value = 105.6693281969054
$ws.Cells.Item(358, 1).Value = 105.9423512847164
$ws.Cells.Item(359, 1).Value = 106.1774153044036
$ws.Cells.Item(360, 1).Value = 106.0802150252765
$ws.Cells.Item(361, 1).Value = 105.845628660529
$ws.Cells.Item(362, 1).Value = 105.6402871043724
$ws.Cells.Item(363, 1).Value = 105.3968974100845
$ws.Cells.Item(364, 1).Value = 106.2972886860027
$ws.Cells.Item(365, 1).Value = 106.4932876351298
$ws.Cells.Item(366, 1).Value = 106.5339593762854
$ws.Cells.Item(367, 1).Value = 106.8306589248023
$ws.Cells.Item(368, 1).Value = 107.0030393253678
$ws.Cells.Item(369, 1).Value = 106.7793415789402
$ws.Cells.Item(370, 1).Value = 106.795781051575
$ws.Cells.Item(371, 1).Value = 106.7500007997975
$ws.Cells.Item(372, 1).Value = 107.029237909112
$ws.Cells.Item(373, 1).Value = 106.9052834268502
$ws.Cells.Item(374, 1).Value = 106.9052834268502
$ws.Cells.Item(375, 1).Value = 106.2139497459497
$ws.Cells.Item(376, 1).Value = 106.3803997553162
$ws.Cells.Item(377, 1).Value = 106.3499104363912
$ws.Cells.Item(378, 1).Value = 105.9045173184649
$ws.Cells.Item(379, 1).Value = 106.0923294302078
$ws.Cells.Item(380, 1).Value = 106.4921973270038
$ws.Cells.Item(381, 1).Value = 106.3354083605602
$ws.Cells.Item(382, 1).Value = 106.111503626499
$ws.Cells.Item(383, 1).Value = 105.8964190156262
$ws.Cells.Item(384, 1).Value = 105.4869723705931
$ws.Cells.Item(385, 1).Value = 105.4017078704993
$ws.Cells.Item(386, 1).Value = 105.4695871550434
$ws.Cells.Item(387, 1).Value = 105.004638655777
$ws.Cells.Item(388, 1).Value = 105.3459263062639
$ws.Cells.Item(389, 1).Value = 105.3728391306908
$ws.Cells.Item(390, 1).Value = 105.3335168646915
$ws.Cells.Item(391, 1).Value = 105.2609090576987
$ws.Cells.Item(392, 1).Value = 105.5975231118028
$ws.Cells.Item(393, 1).Value = 105.8112025071875
$ws.Cells.Item(394, 1).Value = 106.1605772114399
$ws.Cells.Item(395, 1).Value = 106.3060850408301
$ws.Cells.Item(396, 1).Value = 106.1228379107933
$ws.Cells.Item(397, 1).Value = 106.67560884195
$ws.Cells.Item(398, 1).Value = 106.3693415628215
$ws.Cells.Item(399, 1).Value = 106.3747612778059
$ws.Cells.Item(400, 1).Value = 106.4613746140899
$ws.Cells.Item(401, 1).Value = 106.1649106224701
$ws.Cells.Item(402, 1).Value = 106.0595186961931
$ws.Cells.Item(403, 1).Value = 105.7091509971843
$ws.Cells.Item(404, 1).Value = 105.5138861500619
$ws.Cells.Item(405, 1).Value = 104.9399098947893
$ws.Cells.Item(406, 1).Value = 104.8555153658983
$ws.Cells.Item(407, 1).Value = 104.9692890168586
$ws.Cells.Item(408, 1).Value = 104.7840255643677
$ws.Cells.Item(409, 1).Value = 105.1391467272024
$ws.Cells.Item(410, 1).Value = 105.1519289111632
$ws.Cells.Item(411, 1).Value = 105.9757635635999
$ws.Cells.Item(412, 1).Value = 106.0477444053147
$ws.Cells.Item(413, 1).Value = 106.4949211021365
$ws.Cells.Item(414, 1).Value = 106.4689529848464
$ws.Cells.Item(415, 1).Value = 106.4921318860509
$ws.Cells.Item(416, 1).Value = 105.9899404679363
$ws.Cells.Item(417, 1).Value = 106.2597262880127
$ws.Cells.Item(418, 1).Value = 106.0695556738262
$ws.Cells.Item(419, 1).Value = 106.3137715508449
$ws.Cells.Item(420, 1).Value = 106.2200225976687
$ws.Cells.Item(421, 1).Value = 105.645012027633
$ws.Cells.Item(422, 1).Value = 105.8916312355757
$ws.Cells.Item(423, 1).Value = 105.614786853164
$ws.Cells.Item(424, 1).Value = 105.3676607928961
$ws.Cells.Item(425, 1).Value = 105.2252438975294
$ws.Cells.Item(426, 1).Value = 105.5242521061258
$ws.Cells.Item(427, 1).Value = 105.8931432573163
$ws.Cells.Item(428, 1).Value = 105.5991842592717
$ws.Cells.Item(429, 1).Value = 104.9062310556047
$ws.Cells.Item(430, 1).Value = 104.8803251799854
$ws.Cells.Item(431, 1).Value = 104.7230501099798
$ws.Cells.Item(432, 1).Value = 104.5897427196746
$ws.Cells.Item(433, 1).Value = 104.4211044665446
$ws.Cells.Item(434, 1).Value = 104.3818651859989
$ws.Cells.Item(435, 1).Value = 105.3158868947306
$ws.Cells.Item(436, 1).Value = 105.3060571302707
$ws.Cells.Item(437, 1).Value = 105.3143918537409
$ws.Cells.Item(438, 1).Value = 105.43608797624
$ws.Cells.Item(439, 1).Value = 105.44764878235
$ws.Cells.Item(440, 1).Value = 104.9052156448679
$ws.Cells.Item(441, 1).Value = 104.9622494181304
$ws.Cells.Item(442, 1).Value = 105.2352416444732
$ws.Cells.Item(443, 1).Value = 105.3975644643627
$ws.Cells.Item(444, 1).Value = 104.6295540447689
$ws.Cells.Item(445, 1).Value = 104.9808837373804
$ws.Cells.Item(446, 1).Value = 104.7903265781277
$ws.Cells.Item(447, 1).Value = 104.889414250299
$ws.Cells.Item(448, 1).Value = 105.4131064855691
$ws.Cells.Item(449, 1).Value = 105.4418897360076
$ws.Cells.Item(450, 1).Value = 105.4639414387708
$ws.Cells.Item(451, 1).Value = 105.4319629816997
$ws.Cells.Item(452, 1).Value = 104.7936164685284
$ws.Cells.Item(453, 1).Value = 104.145239624369
$ws.Cells.Item(454, 1).Value = 104.3318871102496
$ws.Cells.Item(455, 1).Value = 104.4905410911198
$ws.Cells.Item(456, 1).Value = 103.6085918354658
$ws.Cells.Item(457, 1).Value = 103.7562421601474
$ws.Cells.Item(458, 1).Value = 103.8903760046685
$ws.Cells.Item(459, 1).Value = 103.7859780610992
$ws.Cells.Item(460, 1).Value = 103.5741113403233
$ws.Cells.Item(461, 1).Value = 103.7145324718568
$ws.Cells.Item(462, 1).Value = 104.0842493754243
$ws.Cells.Item(463, 1).Value = 104.3404564729394
$ws.Cells.Item(464, 1).Value = 103.585221735865
$ws.Cells.Item(465, 1).Value = 103.7653157556537
$ws.Cells.Item(466, 1).Value = 104.1210882778258
$ws.Cells.Item(467, 1).Value = 104.631464847381
$ws.Cells.Item(468, 1).Value = 103.7337100117146
$ws.Cells.Item(469, 1).Value = 103.7966996051654
$ws.Cells.Item(470, 1).Value = 103.4101171948362
$ws.Cells.Item(471, 1).Value = 103.6300240784192
$ws.Cells.Item(472, 1).Value = 103.5478577121632
$ws.Cells.Item(473, 1).Value = 102.9486609275069
$ws.Cells.Item(474, 1).Value = 102.5750433996069
$ws.Cells.Item(475, 1).Value = 102.4602606141922
$ws.Cells.Item(476, 1).Value = 103.1788423294094
$ws.Cells.Item(477, 1).Value = 102.9622357049595
$ws.Cells.Item(478, 1).Value = 102.5199102116692
$ws.Cells.Item(479, 1).Value = 102.4718139330352
$ws.Cells.Item(480, 1).Value = 102.8085566656276
$ws.Cells.Item(481, 1).Value = 102.4304232077489
$ws.Cells.Item(482, 1).Value = 102.3094484772367
$ws.Cells.Item(483, 1).Value = 102.4899066534331
$ws.Cells.Item(484, 1).Value = 102.1234883562902
$ws.Cells.Item(485, 1).Value = 102.1801273377591
$ws.Cells.Item(486, 1).Value = 103.3436503791195
$ws.Cells.Item(487, 1).Value = 103.5615093451823
$ws.Cells.Item(488, 1).Value = 103.5541931724525
$ws.Cells.Item(489, 1).Value = 102.9236172810592
$ws.Cells.Item(490, 1).Value = 103.2948919529735
$ws.Cells.Item(491, 1).Value = 103.2463970891344
$ws.Cells.Item(492, 1).Value = 103.3733424720167
$ws.Cells.Item(493, 1).Value = 103.6233704067656
$ws.Cells.Item(494, 1).Value = 103.2697503392256
$ws.Cells.Item(495, 1).Value = 102.793779102072
$ws.Cells.Item(496, 1).Value = 102.6489356307652
$ws.Cells.Item(497, 1).Value = 102.6163238559652
$ws.Cells.Item(498, 1).Value = 102.6163238559652
$ws.Cells.Item(499, 1).Value = 102.681830724076
$ws.Cells.Item(500, 1).Value = 102.4649539071122
$ws.Cells.Item(501, 1).Value = 102.855067107007
$ws.Cells.Item(502, 1).Value = 102.9293400678583
$ws.Cells.Item(503, 1).Value = 103.1041768987263
$ws.Cells.Item(504, 1).Value = 102.9427776461318
$ws.Cells.Item(505, 1).Value = 103.4840684011441
$ws.Cells.Item(506, 1).Value = 103.5212043834226
$ws.Cells.Item(507, 1).Value = 103.1504180895116
$ws.Cells.Item(508, 1).Value = 103.3719954463108
$ws.Cells.Item(509, 1).Value = 103.4370642674428
$ws.Cells.Item(510, 1).Value = 103.7763484279483
$ws.Cells.Item(511, 1).Value = 103.6228990603204
$ws.Cells.Item(512, 1).Value = 103.5738566952796
$ws.Cells.Item(513, 1).Value = 103.4347882192139
$ws.Cells.Item(514, 1).Value = 103.5469482869555
$ws.Cells.Item(515, 1).Value = 103.2467264366051
$ws.Cells.Item(516, 1).Value = 103.1544545203951
$ws.Cells.Item(517, 1).Value = 103.8692719701024
$ws.Cells.Item(518, 1).Value = 104.3842277396642
$ws.Cells.Item(519, 1).Value = 104.3915559864237
$ws.Cells.Item(520, 1).Value = 104.6674623844359
$ws.Cells.Item(521, 1).Value = 104.6898178878533
$ws.Cells.Item(522, 1).Value = 104.7551933689303
$ws.Cells.Item(523, 1).Value = 104.5960437467148
$ws.Cells.Item(524, 1).Value = 104.6902006878935
$ws.Cells.Item(525, 1).Value = 104.7852074684774
$ws.Cells.Item(526, 1).Value = 104.9123355822407
$ws.Cells.Item(527, 1).Value = 105.0186902806776
$ws.Cells.Item(528, 1).Value = 105.1144325244773
$ws.Cells.Item(529, 1).Value = 105.2816393135965
$ws.Cells.Item(530, 1).Value = 105.0828721908445
$ws.Cells.Item(531, 1).Value = 104.8651626675432
$ws.Cells.Item(532, 1).Value = 104.6054256089763
$ws.Cells.Item(533, 1).Value = 104.5663624160557
$ws.Cells.Item(534, 1).Value = 104.7364284131687
$ws.Cells.Item(535, 1).Value = 104.6587141098294
$ws.Cells.Item(536, 1).Value = 105.0745178587069
$ws.Cells.Item(537, 1).Value = 105.3138125047458
$ws.Cells.Item(538, 1).Value = 105.0976556198009
$ws.Cells.Item(539, 1).Value = 105.1056224810769
$ws.Cells.Item(540, 1).Value = 105.0816469881556
$ws.Cells.Item(541, 1).Value = 105.4205196903225
$ws.Cells.Item(542, 1).Value = 105.6311230500221
$ws.Cells.Item(543, 1).Value = 106.0450680769305
$ws.Cells.Item(544, 1).Value = 106.1325401076594
$ws.Cells.Item(545, 1).Value = 106.5395302437839
$ws.Cells.Item(546, 1).Value = 106.4972665330824
$ws.Cells.Item(547, 1).Value = 106.4720283027346
$ws.Cells.Item(548, 1).Value = 106.4853271159013
$ws.Cells.Item(549, 1).Value = 106.3680816024978
$ws.Cells.Item(550, 1).Value = 106.0866152942178
$ws.Cells.Item(551, 1).Value = 106.1587118825117
$ws.Cells.Item(552, 1).Value = 106.8299535068916
$ws.Cells.Item(553, 1).Value = 106.6028286842532
$ws.Cells.Item(554, 1).Value = 106.0814598062294
$ws.Cells.Item(555, 1).Value = 106.1808521381984
$ws.Cells.Item(556, 1).Value = 105.7978409501858
$ws.Cells.Item(557, 1).Value = 105.593380794501
$ws.Cells.Item(558, 1).Value = 105.9407790117796
$ws.Cells.Item(559, 1).Value = 105.7653642339935
$ws.Cells.Item(560, 1).Value = 105.9968193008004
$ws.Cells.Item(561, 1).Value = 106.0344456401505
$ws.Cells.Item(562, 1).Value = 106.1420556762024
$ws.Cells.Item(563, 1).Value = 105.9928747614389
$ws.Cells.Item(564, 1).Value = 106.0271338901329
$ws.Cells.Item(565, 1).Value = 105.9633683256914
$ws.Cells.Item(566, 1).Value = 106.0954421337126
$ws.Cells.Item(567, 1).Value = 106.5053190959491
$ws.Cells.Item(568, 1).Value = 106.8961873605697
$ws.Cells.Item(569, 1).Value = 106.9640958561971
$ws.Cells.Item(570, 1).Value = 107.0336377948922
$ws.Cells.Item(571, 1).Value = 107.4064914732796
$ws.Cells.Item(572, 1).Value = 107.1776948000774
$ws.Cells.Item(573, 1).Value = 107.004986501986
$ws.Cells.Item(574, 1).Value = 106.345505407652
$ws.Cells.Item(575, 1).Value = 106.7004716217405
$ws.Cells.Item(576, 1).Value = 107.0420769701023
$ws.Cells.Item(577, 1).Value = 107.2320711509599
$ws.Cells.Item(578, 1).Value = 107.1437079343426
$ws.Cells.Item(579, 1).Value = 107.3917174856409
$ws.Cells.Item(580, 1).Value = 107.4050874348074
$ws.Cells.Item(581, 1).Value = 107.7080390018057
$ws.Cells.Item(582, 1).Value = 107.5408625768489
$ws.Cells.Item(583, 1).Value = 107.4319561771191
$ws.Cells.Item(584, 1).Value = 107.7441759191631
$ws.Cells.Item(585, 1).Value = 107.7887356439725
$ws.Cells.Item(586, 1).Value = 107.9471756638545
$ws.Cells.Item(587, 1).Value = 108.1152425783343
$ws.Cells.Item(588, 1).Value = 108.1395013012512
$ws.Cells.Item(589, 1).Value = 107.8574282792202
$ws.Cells.Item(590, 1).Value = 108.0929853138811
$ws.Cells.Item(591, 1).Value = 108.2250605160304
$ws.Cells.Item(592, 1).Value = 108.6248811801001
$ws.Cells.Item(593, 1).Value = 108.3762364542303
$ws.Cells.Item(594, 1).Value = 108.6843396646681
$ws.Cells.Item(595, 1).Value = 108.7067199468099
$ws.Cells.Item(596, 1).Value = 108.5033753978106
$ws.Cells.Item(597, 1).Value = 108.4126752314398
$ws.Cells.Item(598, 1).Value = 108.1094936468857
$ws.Cells.Item(599, 1).Value = 108.3882494685005
$ws.Cells.Item(600, 1).Value = 108.528743326328
$ws.Cells.Item(601, 1).Value = 108.1295711474822
$ws.Cells.Item(602, 1).Value = 108.2340513219696
$ws.Cells.Item(603, 1).Value = 108.2630236910356
$ws.Cells.Item(604, 1).Value = 108.2092162607781
$ws.Cells.Item(605, 1).Value = 108.2431945516811
$ws.Cells.Item(606, 1).Value = 108.2310433583919
$ws.Cells.Item(607, 1).Value = 108.4033225662931
$ws.Cells.Item(608, 1).Value = 107.9515422643456
$ws.Cells.Item(609, 1).Value = 108.1333597900188
$ws.Cells.Item(610, 1).Value = 107.7672933387392
$ws.Cells.Item(611, 1).Value = 107.549516610718
$ws.Cells.Item(612, 1).Value = 107.4018191219098
$ws.Cells.Item(613, 1).Value = 107.4010101227956
$ws.Cells.Item(614, 1).Value = 107.8987158763151
$ws.Cells.Item(615, 1).Value = 107.5908649940581
$ws.Cells.Item(616, 1).Value = 107.6618287047726
$ws.Cells.Item(617, 1).Value = 107.6474880562353
$ws.Cells.Item(618, 1).Value = 107.5686421138009
$ws.Cells.Item(619, 1).Value = 107.479401857597
$ws.Cells.Item(620, 1).Value = 107.4666874415297
$ws.Cells.Item(621, 1).Value = 107.6140011476062
